$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update to next match in the week (values shift up; odds refreshed)
$ws.Range("A2").Value = "8UHmSCuQ"
$ws.Range("C2").Value = "22:30"
$ws.Range("E2").Value = "Santa Fe"
$ws.Range("F2").Value = "Chico"
$ws.Range("G2").Value = 1.53
$ws.Range("H2").Value = 3.9
$ws.Range("I2").Value = 6.5
$ws.Range("J2").Value = 2.2
$ws.Range("L2").Value = 7
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("Q2").Value = 2.25
$ws.Range("R2").Value = 1.62
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 2.38
$ws.Range("V2").Value = 1.53
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 6
$ws.Range("Z2").Value = 10
$ws.Range("AB2").Value = 41
$ws.Range("AD2").Value = 8
$ws.Range("AE2").Value = 26
$ws.Range("AF2").Value = 101
$ws.Range("AG2").Value = 201
$ws.Range("AH2").Value = 12
$ws.Range("AI2").Value = 29
$ws.Range("AJ2").Value = 21
$ws.Range("AK2").Value = 81
$ws.Range("AL2").Value = 51
$ws.Range("AM2").Value = 67
$ws.Range("AN2").Value = 3.25
$ws.Range("AO2").Value = 8
$ws.Range("AQ2").Value = 26
$ws.Range("AT2").Value = 2.5
$ws.Range("AU2").Value = 10
$ws.Range("AV2").Value = 81
$ws.Range("AW2").Value = 8
$ws.Range("AX2").Value = 41
$ws.Range("AZ2").Value = 151
$ws.Range("BA2").Value = 201
$ws.Range("BB2").Value = 501

# Row 3: update to next match in the week (values shift up; odds refreshed)
$ws.Range("A3").Value = "hCptA7hl"
$ws.Range("C3").Value = "22:00"
$ws.Range("D3").Value = "MEXICO - LIGA DE EXPANSION MX"
$ws.Range("E3").Value = "Correcaminos"
$ws.Range("F3").Value = "Atl. Morelia"
$ws.Range("G3").Value = 2.8
$ws.Range("H3").Value = 3.05
$ws.Range("I3").Value = 2.47
$ws.Range("J3").Value = 3.35
$ws.Range("K3").Value = 2.07
$ws.Range("L3").Value = 2.95
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 1.31
$ws.Range("P3").Value = 2.9
$ws.Range("Q3").Value = 1.95
$ws.Range("R3").Value = 1.75
$ws.Range("S3").Value = 1.42
$ws.Range("T3").Value = 2.47
$ws.Range("U3").Value = 1.7
$ws.Range("V3").Value = 1.93
$ws.Range("W3").Value = 8.5
$ws.Range("X3").Value = 14
$ws.Range("Y3").Value = 10.25
$ws.Range("Z3").Value = 35
$ws.Range("AA3").Value = 25
$ws.Range("AB3").Value = 32
$ws.Range("AC3").Value = 9
$ws.Range("AD3").Value = 6
$ws.Range("AE3").Value = 13
$ws.Range("AF3").Value = 60
$ws.Range("AG3").Value = 450
$ws.Range("AH3").Value = 8.25
$ws.Range("AI3").Value = 13
$ws.Range("AJ3").Value = 9.25
$ws.Range("AK3").Value = 28
$ws.Range("AL3").Value = 20
$ws.Range("AM3").Value = 28
$ws.Range("AN3").Value = 4.75
$ws.Range("AO3").Value = 15
$ws.Range("AP3").Value = 21
$ws.Range("AQ3").Value = 70
$ws.Range("AR3").Value = 100
$ws.Range("AS3").Value = 250
$ws.Range("AT3").Value = 2.6
$ws.Range("AU3").Value = 6.4
$ws.Range("AV3").Value = 50
$ws.Range("AW3").Value = 4.45
$ws.Range("AX3").Value = 12.5
$ws.Range("AY3").Value = 18
$ws.Range("AZ3").Value = 50
$ws.Range("BA3").Value = 70
$ws.Range("BB3").Value = 200
$ws.Range("BC3").Value = 51
$ws.Range("BD3").Value = 51

# Row 4 no longer exists once data shifts up; remove it and its dimension entry
$ws.Rows.Item(4).Delete()
